$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.507.44"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").Value = "1.671.75"
$ws.Range("E3").Value = "  +1.64%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5285"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.24%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2679"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06373"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07805"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.87%  "
$ws.Range("D12").Value = "1.677.44"
$ws.Range("E12").Value = "  +1.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.484"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5564"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.16%  "
$ws.Range("D17").Value = "26.498.78"
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.764"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.304"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1269"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.378"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.427"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.82%  "
$ws.Range("E29").Value = "  +4.76%  "
$ws.Range("E30").Value = "  +1.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.606"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.99%  "
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.686"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.007"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("E35").Value = "  +9.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.422"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.784"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.041"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.28%  "
$ws.Range("D40").Value = "1.091.57"
$ws.Range("E40").Value = "  +6.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8565"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.11%  "
$ws.Range("D44").Value = "1.816.18"
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "58.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.78%  "
$ws.Range("D46").Value = "0.0₈107"
$ws.Range("E46").Value = "  -2.20%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.002"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.114"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.517"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.05%  "
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.005"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.58%  "
